# Applies the PlayerPerformance_6547 update:
#  1. Clear the stray empty INNING_NUMBER cells (B2/B8/B14) on "ODI Batting"
#     so they disappear entirely instead of round-tripping as empty cells.
#  2. Normalise the non-breaking space in "ODI Batting"!E14 to a plain space.
#  3. Append a new "ODI Batting Extra" worksheet (after "ODI Bowling") with
#     MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL /
#     MAN_OF_MATCH columns for 19 matches.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: fix up "ODI Batting" -------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()
$batting.Range("B8").ClearContents()
$batting.Range("B14").ClearContents()
$batting.Range("E14").Value = " "

# --- 3: add the "ODI Batting Extra" sheet ----------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowling)
$extra.Name = "ODI Batting Extra"

# Match the page margins used by the sheets already in the workbook.
$extra.PageSetup.LeftMargin = 0.75 * 72
$extra.PageSetup.RightMargin = 0.75 * 72
$extra.PageSetup.TopMargin = 1 * 72
$extra.PageSetup.BottomMargin = 1 * 72
$extra.PageSetup.HeaderMargin = 0.5 * 72
$extra.PageSetup.FooterMargin = 0.5 * 72

# MATCH_CODE, NUM_4, NUM_6 and PERCENT_RUNS_OF_TOTAL are id-like / percentage-text
# columns (e.g. "41.64%") that must stay text rather than being auto-converted to
# numbers, so format those columns as Text up front.
$extra.Range("A2:A20").NumberFormat = "@"
$extra.Range("C2:E20").NumberFormat = "@"

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $extra.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Match the bold/centred/bordered header style already used by the other sheets
# (re-using the existing style avoids minting a near-duplicate one).
$bowling.Range("A1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
    @("4563", 6,    $null, $null, $null,    "NO"),
    @("4566", 6,    "0",   "0",   "0.38%",  "NO"),
    @("4568", $null, $null, $null, $null,   "NO"),
    @("4605", 7,    "10",  "7",   "41.64%", "YES"),
    @("4608", 7,    "3",   "3",   "19.18%", "YES"),
    @("4614", 7,    "1",   "1",   "5.83%",  "NO"),
    @("4625", 7,    $null, $null, $null,    "NO"),
    @("4636", 7,    "5",   "0",   "16.32%", "NO"),
    @("4639", 6,    "0",   "0",   "2.83%",  "NO"),
    @("4642", $null, $null, $null, $null,   "NO"),
    @("4647", $null, $null, $null, $null,   "NO"),
    @("4648", 6,    "1",   "0",   "14.63%", "NO"),
    @("4673", $null, $null, $null, $null,   "NO"),
    @("4686", $null, $null, $null, $null,   "NO"),
    @("4688", 7,    "0",   "0",   "3.07%",  "NO"),
    @("4690", $null, $null, $null, $null,   "NO"),
    @("4692", $null, $null, $null, $null,   "NO"),
    @("4695", 7,    "4",   "0",   "20.37%", "NO"),
    @("4697", 7,    "3",   "1",   "8.81%",  "NO")
)

$r = 2
foreach ($row in $rows) {
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($null -ne $row[1]) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        $extra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($null -ne $row[3]) {
        $extra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($null -ne $row[4]) {
        $extra.Cells.Item($r, 5).Value = $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Restore the original active sheet/selection (the edit shouldn't change which
# tab is shown when the workbook is reopened).
$wb.Worksheets.Item("Player Info").Activate()
